# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "last updated" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 28 de Abril de 2020 a las 22:52"

# 2. Update Estados Unidos totals (row 4)
$ws.Range("B4").Value = 1029179
$ws.Range("C4").Value = 18823
$ws.Range("E4").Value = 830473
$ws.Range("F4").Value = 14868
$ws.Range("G4").Value = 1771
$ws.Range("H4").Value = 58568

# 3. Insert Costa Rica ahead of Honduras (alphabetical-ish reorder), pushing
#    Honduras and Niger down a row, and give Costa Rica its updated numbers.
# Row 100 was Honduras -> becomes Costa Rica with new figures
$ws.Range("A100").Value = "Costa Rica"
$ws.Range("B100").Value = 705
$ws.Range("C100").Value = 8
$ws.Range("D100").Value = 306
$ws.Range("E100").Value = 393
$ws.Range("F100").Value = 8
$ws.Range("G100").Value = 0
$ws.Range("H100").Value = 6

# Row 101 was Niger -> becomes Honduras (carrying Honduras' old figures)
$ws.Range("A101").Value = "Honduras"
$ws.Range("B101").Value = 702
$ws.Range("C101").Value = 41
$ws.Range("D101").Value = 79
$ws.Range("E101").Value = 559
$ws.Range("F101").Value = 10
$ws.Range("G101").Value = 3
$ws.Range("H101").Value = 64

# Row 102 was Costa Rica -> becomes Niger (carrying Niger's old figures)
$ws.Range("A102").Value = "Niger"
$ws.Range("B102").Value = 701
$ws.Range("C102").Value = 0
$ws.Range("D102").Value = 385
$ws.Range("E102").Value = 287
$ws.Range("F102").Value = 0
$ws.Range("G102").Value = 0
$ws.Range("H102").Value = 29

# Row 103 (Burkina Faso) is unchanged.

# 4. Update San Marino totals (row 106)
$ws.Range("B106").Value = 553
$ws.Range("C106").Value = 15
$ws.Range("E106").Value = 448
$ws.Range("F106").Value = 5
